# Actualización desde MV -datos-
# Add a new "Agosto.2021" period column (BH), mirroring the pattern already
# used in the sheet: each new reporting period column is added to the right
# of the previous last column (BG) and carries forward that column's values
# (the data hasn't changed since the last bulletin yet, so the new period
# repeats the prior one) along with the same header formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in BH1, with the same value + style as the rest of the header row
$ws.Range("BH1").Value = "Agosto.2021"
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Copy BG2:BG19 values into BH2:BH19 (same forward-filled pattern as BF->BG)
for ($r = 2; $r -le 19; $r++) {
    $src = $ws.Cells.Item($r, 59)   # column BG
    $dst = $ws.Cells.Item($r, 60)   # column BH
    $dst.Value = $src.Value2
}
